$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column widths (Excel snaps width to whole-pixel increments, so we pick the
# ColumnWidth input that rounds to the closest achievable raw width)
$ws.Range("E1").ColumnWidth = 13.86
$ws.Range("G1").ColumnWidth = 9.14

# Row 2
$ws.Range("C2").Value = 490
$ws.Range("D2").Value = 94443.533905760924
$ws.Range("F2").Value = 0.99483849405114932

# Row 3
$ws.Range("C3").Value = 490
$ws.Range("D3").Value = 10.471349452636249
$ws.Range("E3").Value = 0.0012944081158074994
$ws.Range("F3").Value = 0.020922974839795976
$ws.Range("H3").Value = "**"

# Row 4
$ws.Range("C4").Value = 490
$ws.Range("D4").Value = 900.1929335880933
$ws.Range("F4").Value = 0.64753093749706525
$ws.Range("H4").Value = "***"

# Row 5
$ws.Range("C5").Value = 490
$ws.Range("D5").Value = 147.36831823318715
$ws.Range("F5").Value = 0.23121374881277212
$ws.Range("H5").Value = "***"

# G3 text change: "small to medium" -> "small"
$ws.Range("G3").Value = "small"
